
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert one new observation as the new row 378.
# This shifts every following row (old 378..406) down by one (to 379..407),
# preserving their original data, and populates the new row 378 with
# the latest weekly data point.
$ws.Rows.Item(378).Insert()

$ws.Cells.Item(378, 1).Value = 9
$ws.Cells.Item(378, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(378, 3).Value = "Metropolitana"
$ws.Cells.Item(378, 4).Value = 45013
$ws.Cells.Item(378, 5).Value = 13
$ws.Cells.Item(378, 6).Value = 300000001
$ws.Cells.Item(378, 7).Value = "Rabanito"
$ws.Cells.Item(378, 8).Value = "Sin especificar"
$ws.Cells.Item(378, 9).Value = "Primera"
$ws.Cells.Item(378, 10).Value = 7000
$ws.Cells.Item(378, 11).Value = 3000
$ws.Cells.Item(378, 12).Value = 3000
$ws.Cells.Item(378, 13).Value = 3000
$ws.Cells.Item(378, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(378, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(378, 16).Value = 30
$ws.Cells.Item(378, 17).Value = 100
$ws.Cells.Item(378, 18).Value = "Hortaliza"

$ws.Cells.Item(378, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
